# ReportTemplate.xlsx edit: rename sheet, restyle header/body rows,
# add a border row, resize column H, update one placeholder string,
# and trim two trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Rename the sheet
# ---------------------------------------------------------------
$ws.Name = "Registros"

# ---------------------------------------------------------------
# 2. Update the placeholder text in H3 (adds a new shared string)
# ---------------------------------------------------------------
$ws.Range("H3").Value = "{d.dados[+i].obs}"

# ---------------------------------------------------------------
# 3. Style the header row (A1:H1): bold 14pt Calibri, light-blue
#    fill, thin black box border.
# ---------------------------------------------------------------
$header = $ws.Range("A1:H1")
$header.Font.Name = "Calibri"
$header.Font.Size = 14
$header.Font.Bold = $true
$header.Font.ThemeColor = 1
$header.Interior.Color = 16308937
$header.Interior.PatternColor = 16308937
$header.Borders.LineStyle = 1
$header.Borders.Color = 0

# ---------------------------------------------------------------
# 4. Style row 2 / row 3 body cells with the new Calibri 11 look.
# ---------------------------------------------------------------

# A2, H2, A3, B3 -> plain Calibri 11, automatic (theme) color
# (the COM bridge does not fan a property write out across a
# multi-area union, so each cell is addressed individually)
foreach ($addr in @("A2", "H2", "A3", "B3")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.ThemeColor = 1
}

# B2 -> Calibri 11, explicit black color
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Color = 0
$ws.Range("B2").WrapText = $false

# C2:G2, C3:G3 -> Calibri 11, explicit black color, right aligned
foreach ($addr in @("C2:G2", "C3:G3")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.Color = 0
    $c.HorizontalAlignment = -4152
}

# H3 -> Calibri 11, no explicit color
$ws.Range("H3").Font.Name = "Calibri"
$ws.Range("H3").Font.Size = 11

# ---------------------------------------------------------------
# 5. Row 4: blank cells with a thin top border (Arial, theme color)
# ---------------------------------------------------------------
$row4 = $ws.Range("A4:H4")
$row4.Font.Name = "Arial"
$row4.Font.ThemeColor = 1
$row4.Borders.Item(8).LineStyle = 1
$row4.Borders.Item(8).Color = 0

# ---------------------------------------------------------------
# 6. Widen column H
# ---------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 28.6

# ---------------------------------------------------------------
# 7. Trim the two trailing blank rows (999, 1000)
# ---------------------------------------------------------------
$ws.Rows.Item(1000).Delete()
$ws.Rows.Item(999).Delete()
